$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 75 (shifts existing rows 75..152 down to 76..153)
$ws.Rows.Item(75).EntireRow.Insert()

# Populate the new row 75 with the new Papaya price record
$ws.Cells.Item(75, 1).Value = 10
$ws.Cells.Item(75, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(75, 3).Value = "La Araucanía"
$ws.Cells.Item(75, 4).Value = 45280
$ws.Cells.Item(75, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(75, 5).Value = 9
$ws.Cells.Item(75, 6).Value = "Fruta"
$ws.Cells.Item(75, 7).Value = 100108
$ws.Cells.Item(75, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(75, 9).Value = 100108004
$ws.Cells.Item(75, 10).Value = "Papaya"
$ws.Cells.Item(75, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(75, 12).Value = "Primera"
$ws.Cells.Item(75, 13).Value = 80
$ws.Cells.Item(75, 14).Value = 24000
$ws.Cells.Item(75, 15).Value = 24000
$ws.Cells.Item(75, 16).Value = 24000
$ws.Cells.Item(75, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(75, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(75, 19).Value = 2400
$ws.Cells.Item(75, 20).Value = 10
